$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "Alessio Zandonai"
$ws.Range("B30").Value = "Daniele  Dalbosco | iMontagna"
$ws.Range("C30").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("D30").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E30").Value = "Luca Frasca | Clitoriders"
$ws.Range("F30").Value = "Davide  Bazzano  | iMontagna"
